$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update pt_max column (E) values from 50 to 70 for rows 2-7
$ws.Range("E2:E7").Value = 70
